$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prepare formatting for the new rows (23-26) by copying the blank,
# already-formatted row 22 down. This gives every new cell the same
# bordered / centered base style ("s=4") that the rest of the table uses,
# matching what a user would get by dragging/copying the last template row.
$ws.Range("A22:M22").Copy()
$ws.Range("A23:M26").PasteSpecial(-4122)  # xlPasteFormats

# Turn on word-wrap for the "Copper coins" name cell ahead of time (before
# any values are typed in) so the new cell-format record is created in the
# same order as in the saved workbook, without disturbing shared-string order.
$ws.Range("F23").WrapText = $true

# --- Row 22: OML token item ---
$ws.Range("A22").Value = "'0001"
$ws.Range("C22").Value = "'0001"
$ws.Range("D22").Value = "OML"
$ws.Range("E22").Value = "OML"
$ws.Range("F22").Value = "OML"
$ws.Range("G22").Value = "OML代币，可在OML游戏中通过不同途径获得"
$ws.Range("H22").Value = "OML代幣，可在OML遊戲中透過不同途徑獲得"
$ws.Range("I22").Value = "OML tokens can be obtained through different ways in OML games"
$ws.Range("L22").Value = 1
$ws.Range("M22").Value = 0

# --- Row 23: Copper coins item ---
$ws.Range("A23").Value = "'0002"
$ws.Range("C23").Value = "'0002"
$ws.Range("D23").Value = "铜钱"
$ws.Range("E23").Value = "銅錢"
$ws.Range("F23").Value = "Copper coins"
$ws.Range("G23").Value = "OML游戏的通用货币，可以购买游戏中的各种道具"
$ws.Range("H23").Value = "OML遊戲的通用貨幣，可以購買遊戲中的各種道具"
$ws.Range("I23").Value = "The universal currency of OML games, which can be used to purchase various props in the game"
$ws.Range("L23").Value = 1
$ws.Range("M23").Value = 0

# --- Row 24: Stamina/physical strength item ---
$ws.Range("A24").Value = "'0003"
$ws.Range("C24").Value = "'0003"
$ws.Range("D24").Value = "体力"
$ws.Range("E24").Value = "體力"
$ws.Range("F24").Value = "physical strength"
$ws.Range("G24").Value = "体力每过一小时可以增加一点，上限是24点"
$ws.Range("H24").Value = "體力每過一小時可以增加一點，上限是24點"
$ws.Range("I24").Value = "Physical strength can be increased by one point every hour, with a maximum limit of 24 points."
$ws.Range("L24").Value = 1
$ws.Range("M24").Value = 0

# Rows 25 and 26 stay blank (formatting only, already applied above).

# --- Column widths: split the combined column 9-10 sizing so column 9
# gets its own (wider) width, leaving column 10 unchanged. ---
$ws.Columns.Item(9).ColumnWidth = 104.14

# --- Selection moves to G29 to match the saved view state. ---
$ws.Range("G29").Select()
